# Update column G ("K") on Sheet1 of the Tyler Anderson 2022 save_data
# workbook. The save_data pipeline was regenerated to source strikeouts
# ("K") directly instead of the old "Strike#" value, and std/mean +
# s_vals were recalculated. The net effect captured by the diff is a
# straightforward overwrite of the existing G2:G34 values with the new
# K counts; no other cells, columns, or formatting changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 6
    3  = 10
    4  = 3
    5  = 6
    6  = 2
    7  = 3
    8  = 3
    9  = 4
    10 = 6
    11 = 4
    12 = 3
    13 = 3
    14 = 4
    15 = 6
    16 = 4
    17 = 4
    18 = 6
    19 = 2
    20 = 2
    21 = 8
    22 = 2
    23 = 5
    24 = 6
    25 = 8
    26 = 7
    27 = 5
    28 = 7
    29 = 3
    30 = 4
    31 = 4
    32 = 4
    33 = 3
    34 = 7
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
